# Insert a new data row at row 104 (pushing the existing row 104..190
# down to 105..191) and populate the new row with a fresh record.
# This matches the commit "Fruta / hortaliza, semanal" which adds one
# new weekly observation into the middle of the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("104:104").Insert()

$ws.Range("A104").Value = 5
$ws.Range("B104").Value = "Macroferia Regional de Talca"
$ws.Range("C104").Value = "Maule"
$ws.Range("D104").Value = 44566
$ws.Range("E104").Value = 7
$ws.Range("F104").Value = 100112008
$ws.Range("G104").Value = "Coliflor"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 2000
$ws.Range("K104").Value = 800
$ws.Range("L104").Value = 800
$ws.Range("M104").Value = 800
$ws.Range("N104").Value = "$/unidad"
$ws.Range("O104").Value = "Región del Maule"
$ws.Range("P104").Value = 800
$ws.Range("Q104").Value = 1
$ws.Range("R104").Value = "Hortaliza"
